# Update "Forecast Comparison" sheet with corrected forecast output:
#  - Insert a new column B "Week_Start_Date" holding the first date of each forecast week
#  - Shorten the Week labels in column A from W01..W09 to W1..W9 (W10..W16 stay the same)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before the current column B ("ASIN"), shifting the rest right.
$ws.Columns.Item(2).Insert()

# Header for the newly inserted column.
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# Week labels for column A (leading zero dropped for weeks 1-9).
$weekLabels = @("W1","W2","W3","W4","W5","W6","W7","W8","W9","W10","W11","W12","W13","W14","W15","W16")

# Start date of each forecast week, for the new column B.
$weekStarts = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

for ($i = 0; $i -lt $weekLabels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $weekLabels[$i]
    # Prefix with an apostrophe so the date-looking text is stored as plain text,
    # not auto-converted into a date serial number.
    $ws.Cells.Item($row, 2).Value = "'" + $weekStarts[$i]
}
